# Actualizacion automatica del tracker
# 1) Rellenar resultados pendientes (G/H) de partidos ya finalizados con "Fallo"
# 2) Anadir nuevas filas de pronosticos al final de la tabla

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Completar resultados "Fallo" / -1 para filas pendientes ---
$falloRows = @(123, 124, 125, 135, 137, 138)
foreach ($r in $falloRows) {
    $ws.Cells.Item($r, 7).Value = "Fallo"
    $ws.Cells.Item($r, 8).Value = -1
}

# --- 2) Anadir nuevas filas de pronosticos (141-150) ---
# event_id, fecha, jugador_A, jugador_B, pronostico, cuota
$newRows = @(
    @(14680553, "2025-09-15", "Thiago Monteiro", "Daniel Merida", "Gana Daniel Merida", 1.83),
    @(14684248, "2025-09-15", "Tom Paris", "Maxime Janvier", "Gana Maxime Janvier", 2),
    @(14684202, "2025-09-15", "Florent Bax", "Etienne Donnet", "Gana Florent Bax", 2),
    @(14685798, "2025-09-15", "Alexey Vatutin", "Louis Tessa", "Gana Louis Tessa", 4.5),
    @(14681205, "2025-09-14", "Carlo Alberto Caniato", "Preston Stearns", "Gana Preston Stearns", 3.5),
    @(14681213, "2025-09-14", "Samir Banerjee", "Elmar Ejupovic", "Gana Elmar Ejupovic", 4),
    @(14681211, "2025-09-14", "Daniel Milavsky", "Blaise Bicknell", "Gana Blaise Bicknell", 2.1),
    @(14686098, "2025-09-15", "Gerard Campana Lee", "Matei Varbanciu", "Gana Matei Varbanciu", 13),
    @(14686095, "2025-09-15", "Juan Cruz Martin Manzano", "Thomas Faurel", "Gana Thomas Faurel", 1.83),
    @(14686090, "2025-09-15", "Dan Martin", "Radu Mihai Papoe", "Gana Radu Mihai Papoe", 1.36)
)

$startRow = 141
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    # leading apostrophe keeps the ISO date as literal text (not a date serial)
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    # resultado / profit still pending -> empty text placeholders, same as
    # the other not-yet-settled rows already in the sheet
    $ws.Cells.Item($r, 7).Value = "'"
    $ws.Cells.Item($r, 8).Value = "'"
}
